$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 76

# Set the date-label cell as plain text (it looks like a date, so use a
# leading quote to force text entry, then reset the style back to the
# sheet's default "Normal" so no stray number-format style is left behind).
$ws.Cells.Item($row, 1).Value = "'01-07-2021"
$ws.Cells.Item($row, 1).Style = "Normal"

$values = @(
    -11817,
    -11503,
    0,
    -314,
    7583,
    6273,
    8675,
    1,
    0,
    -2403,
    -754,
    -223,
    -1464,
    248,
    685,
    -818,
    0,
    -717,
    -31,
    -70,
    2883,
    624,
    2316,
    5,
    -62,
    -3192,
    -57,
    -3461,
    41,
    285
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
